$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.443.08"
$ws.Range("D3").Value = "'1.856.12"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'244.83"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'0.6949"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D8").Value = "'0.07677"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "'0.3063"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").Value = "'0.07776"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "'5.149"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "'1.852.99"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "'90.97"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "'6.323"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "'29.426.25"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "'0.000008299"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "'2.098.74"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").Value = "'237.78"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("D23").Value = "'7.619"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").Value = "'160.01"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("D27").Value = "'8.877"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "'18.25"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "'1.527"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "'4.240"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "'4.151"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'1.207"
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("D33").Value = "'0.05107"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'0.7694"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "'2.678"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'1.329.74"
$ws.Range("E38").Value = "  +7.10%  "
$ws.Range("D39").Value = "'0.01870"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").Value = "'0.9511"
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("D42").Value = "'105.99"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "'5.821"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'9.823"
$ws.Range("E45").Value = "  +2.60%  "
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("D47").Value = "'1.998.87"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").Value = "'1.781"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").Value = "'63.28"
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").Value = "'6.967"
$ws.Range("E51").Value = "  +0.74%  "
